$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the expense data rows (3-6): clear both contents and formatting
# so the now-empty rows/cells disappear from the sheet, leaving only
# row 3's own row-level formatting behind (matches "delete_data" behavior).
$ws.Range("A3:C6").Clear()
